$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="29.500.49"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +0.92%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("D3").Formula = '="1.871.60"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("D5").Formula = '="0.7178"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +0.90%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="239.26"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = '="  +0.59%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("E7").Formula = '="  +0.05%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("D8").Formula = '="0.07820"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = '="  -4.52%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("D9").Formula = '="0.3071"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  +0.96%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("D10").Formula = '="25.31"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="  +8.96%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="0.08229"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = '="1.874.41"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  -0.79%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("D13").Formula = '="0.7218"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +1.79%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("E14").Formula = '="  +1.12%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("D15").Formula = '="89.97"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="  +0.32%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D16").Formula = '="29.548.12"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = '="  +1.00%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("D17").Formula = '="5.824"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = '="  +0.58%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("D18").Formula = '="0.000007849"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Formula = '="  -1.21%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="240.94"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  +1.53%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="13.31"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = '="  -0.43%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="2.142.69"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = '="  +0.86%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  -0.02%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)
$ws.Range("E23").Formula = '="  +0.13%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="7.737"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D25").Formula = '="0.1565"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  +7.33%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="162.73"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +0.05%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="8.957"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = '="  -0.08%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("D28").Formula = '="18.30"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = '="  +1.13%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="1.359"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Formula = '="  -4.70%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="4.334"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = '="  -1.34%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("E33").Formula = '="  +1.12%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +0.43%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  +2.32%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="0.7164"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = '="  +1.18%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  -0.03%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="2.672"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Formula = '="  -0.03%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  +0.54%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  -0.38%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="1.173.87"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = '="  +2.81%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="0.9070"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  -1.82%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="5.994"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  +1.64%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="0.4309"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +0.55%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)
$ws.Range("D45").Formula = '="71.47"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  +1.40%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="1.002"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Formula = '="  +0.14%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)
$ws.Range("D47").Formula = '="102.34"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  -0.50%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)
$ws.Range("D48").Formula = '="0.5361"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  -0.56%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)
$ws.Range("D49").Formula = '="1.766"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = '="  -0.60%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)
$ws.Range("D50").Formula = '="9.160"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  -0.60%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("D51").Formula = '="7.019"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = '="  +0.94%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)
